$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect the new test data
# Columns: A=FirstName, B=LastName, C=Email, D=Password, E=ConfirmPassword
$ws.Range("A2").Value = "My"
$ws.Range("B2").Value = "Name"
$ws.Range("C2").Value = "ak1217@gmail.com"
$ws.Range("D2").Value = "Heoo12334"
$ws.Range("E2").Value = "Heoo12334"

# Update the selected cell/range on the sheet view
$ws.Range("E2").Select()
